$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...as there were a few options available going forward" gains a
# new parenthetical citation right before the following full stop:
# "...going forward (ref [1])."
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("as there were a few options available going forward", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(0)
    $r.InsertAfter(" (ref [1])")
}

Write-Host "done change 1"

# ---------------------------------------------------------------------------
# Change 2: insert a new definition of "real-time" plus a survey suggestion
# between "...can be performed in real-time." and "Evaluation will be
# performed by measuring..."
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("Evaluation will be performed by measuring the performance of the object manipulations and", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Collapse(1)
    $insertText = "Real-time will be defined as providing the necessary feedback to the user after specific inputs within a time-frame small enough that it feels instant, roughly less than 100 milliseconds. The project could benefit from a survey ran where participants are asked whether they felt the software performed in real-time. "
    $r2.InsertBefore($insertText)
}

Write-Host "done change 2"

# ---------------------------------------------------------------------------
# Helper functions used to build the new paragraphs appended at the very end
# of the document (References section + reviewer-style red-text notes).
# ---------------------------------------------------------------------------
function Get-LastRealParagraph {
    param($doc)
    $idx = $doc.Paragraphs.Count - 1
    return $doc.Paragraphs.Item($idx)
}

function New-ParaAfter {
    # Inserts a brand-new paragraph right after $prevPara and returns it.
    param($doc, $prevPara, $firstLineIndentPt, $lineSpacingRule)
    $anchor = $prevPara.Range
    $anchor.Collapse(0)
    $anchor.InsertParagraphAfter()
    $p = Get-LastRealParagraph $doc
    if ($lineSpacingRule -ne $null) {
        $p.Format.LineSpacingRule = $lineSpacingRule
    }
    if ($firstLineIndentPt -ne $null) {
        $p.Format.FirstLineIndent = $firstLineIndentPt
    }
    return $p
}

function Add-FormattedRun {
    # Appends $text just before the paragraph mark of $para, applying the
    # requested character formatting, and returns the inserted Range.
    param($doc, $para, $text, $bold, $colorVal, $sizePt, $fontName)
    $endPos = $para.Range.End - 1
    $insertRange = $doc.Range($endPos, $endPos)
    $insertRange.InsertAfter($text)
    $newEndPos = $endPos + $text.Length
    $runRange = $doc.Range($endPos, $newEndPos)
    if ($sizePt -ne $null) { $runRange.Font.Size = $sizePt }
    if ($bold) { $runRange.Font.Bold = $true }
    if ($colorVal -ne $null) { $runRange.Font.Color = $colorVal }
    if ($fontName -ne $null) { $runRange.Font.Name = $fontName }
    return $runRange
}

$RED = 255          # wdColorRed       -> FF0000
$GREEN_ACCENT = 4697456   # theme accent6 70AD47, expressed as a literal RGB

# ---------------------------------------------------------------------------
# Change 4a: the trailing "_GoBack" bookmark used to sit between "...can" and
# " perform simple manipulations...". It now belongs further down, inside the
# new "Possible surveys" note, so drop it from its old spot first.
# ---------------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
}

# ---------------------------------------------------------------------------
# Change 4b: append the new paragraphs after the final "Evaluation" paragraph.
# ---------------------------------------------------------------------------
$lastPara = Get-LastRealParagraph $d

# Empty spacer paragraph.
$pEmpty = New-ParaAfter $d $lastPara 18 1

# "References" heading (bold), no first-line indent.
$pRefHeading = New-ParaAfter $d $pEmpty 0 1
Add-FormattedRun $d $pRefHeading "References" $true $null 12 $null

# "[1] https://developer.microsoft.com/en-us/windows/mixed-reality/development"
$pRefItem = New-ParaAfter $d $pRefHeading 0 1
Add-FormattedRun $d $pRefItem "[1] " $false $null 12 $null
Add-FormattedRun $d $pRefItem "https://developer.microsoft.com/en-us/windows/mixed-reality/development" $false $null 12 $null

# Reviewer note: "Define real-time"
$pNote1 = New-ParaAfter $d $pRefItem 18 1
Add-FormattedRun $d $pNote1 ([char]0x2714) $false $GREEN_ACCENT $null "Segoe UI Symbol"
Add-FormattedRun $d $pNote1 " " $false $null $null "Segoe UI Symbol"
Add-FormattedRun $d $pNote1 "Define" $false $RED 12 $null
Add-FormattedRun $d $pNote1 " real-time" $false $RED 12 $null

# Reviewer note: "Reference research done on software"
$pNote2 = New-ParaAfter $d $pNote1 18 1
Add-FormattedRun $d $pNote2 ([char]0x2714) $false $GREEN_ACCENT $null "Segoe UI Symbol"
Add-FormattedRun $d $pNote2 " " $false $GREEN_ACCENT $null "Segoe UI Symbol"
Add-FormattedRun $d $pNote2 "Ref" $false $RED 12 $null
Add-FormattedRun $d $pNote2 "erence research done on software" $false $RED 12 $null

# Reviewer note: "Possible surveys: are you convinced this has been
# manipulated in front of them?" -- carries the relocated _GoBack bookmark.
$pNote3 = New-ParaAfter $d $pNote2 18 1
Add-FormattedRun $d $pNote3 ([char]0x2714) $false $GREEN_ACCENT $null "Segoe UI Symbol"
Add-FormattedRun $d $pNote3 " " $false $GREEN_ACCENT $null "Segoe UI Symbol"

$bmPos = $pNote3.Range.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Add-FormattedRun $d $pNote3 "Possible surveys" $false $RED 12 $null
Add-FormattedRun $d $pNote3 ": are you convinced this has been manipulated in front of them" $false $RED 12 $null
Add-FormattedRun $d $pNote3 "?" $false $RED 12 $null

Write-Host "done change 4"
